{"js": "// Replace the outgoing signatory's name \"LEONARDO SILVERIO FERREIRA\" with the\n// new signatory's name \"MANOEL JEFETE DA SILVA TENORIO\" in the approval\n// table of the NR-10 OEM document, preserving the run's existing formatting\n// (bold, color, etc.).\n\nconst oldName = \"LEONARDO SILVERIO FERREIRA\";\nconst newName = \"MANOEL JEFETE DA SILVA TENORIO\";\n\nconst body = context.document.body;\nconst results = body.search(oldName, { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  // insertText with Replace keeps the run's existing character formatting\n  // (bold/color) intact while only swapping the visible text.\n  results.items[i].insertText(newName, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the outgoing signatory's name \"LEONARDO SILVERIO FERREIRA\" with the\n# new signatory's name \"MANOEL JEFETE DA SILVA TENORIO\" in the approval\n# table of the NR-10 OEM document.\n\n$d = $word.ActiveDocument\n\n$oldName = \"LEONARDO SILVERIO FERREIRA\"\n$newName = \"MANOEL JEFETE DA SILVA TENORIO\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldName\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newName\n$find.Forward = $true\n$find.Wrap = 1        # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.MatchWildcards = $false\n\n$find.Execute([ref]$oldName, $false, $true, $false, $false, $false, $true, 1, $false, [ref]$newName, 2) | Out-Null\n"}
